# Generate Report for Handoff
# Adds a new tracked file (7c68ea15-415e-4320-b920-21e338c7e652...md) to the
# localization status workbook: one new summary row on "Overview", and one
# new detail row on each locale sheet ("zh-cn", "de-de").

$wb = $excel.ActiveWorkbook

$newFileName   = "7c68ea15-415e-4320-b920-21e338c7e652ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newPathName   = "e2e\7c68ea15-415e-4320-b920-21e338c7e652ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newExtension  = ".md"
$statusText    = "Ready for handoff"
$hoDateTime    = "2016-08-17 18:27:20"
$zhHandoffDate = "2016-08-17 18:27:15"
$deHandoffDate = "2016-08-17 18:27:20"
$zhXliff       = "7c68ea15-415e-4320-b920-21e338c7e652oooooooooooooooooooooooooooooooooooooooo.f11a6c1cecff46dbd1fec1311736849712c32782.zh-cn.xlf"
$deXliff       = "7c68ea15-415e-4320-b920-21e338c7e652oooooooooooooooooooooooooooooooooooooooo.f11a6c1cecff46dbd1fec1311736849712c32782.de-de.xlf"
$commitSha     = "bcfb89c480396aea515646f244fb157714ebb5c4"
$newFileUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newFileName"

# ---------------------------------------------------------------------
# Overview sheet: new row 3
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newFileName
$wsOverview.Range("C3").Value = $newExtension
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Range("G3").Value = $hoDateTime
$wsOverview.Range("G3").NumberFormat = $wsOverview.Range("G2").NumberFormat

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newFileUrl, "", "", $newPathName)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------
# zh-cn sheet: new row 3
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B3").Value = $newExtension
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("F3").Style = "Normal"
$wsZh.Range("G3").Value = $zhXliff
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("H3").NumberFormat = $wsZh.Range("H2").NumberFormat
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = $wsZh.Range("K2").NumberFormat
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("M3").Style = "Normal"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "'False"
$wsZh.Range("O3").Style = "Normal"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newFileUrl, "", "", $newFileName)

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

$wsZh.Columns.Item(3).ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------
# de-de sheet: new row 3
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B3").Value = $newExtension
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("F3").Style = "Normal"
$wsDe.Range("G3").Value = $deXliff
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("H3").NumberFormat = $wsDe.Range("H2").NumberFormat
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = $wsDe.Range("K2").NumberFormat
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("M3").Style = "Normal"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "'False"
$wsDe.Range("O3").Style = "Normal"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newFileUrl, "", "", $newFileName)

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))

$wsDe.Columns.Item(3).ColumnWidth = 17.2159881591797
